$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'30.111.04"
$ws.Range("E2").Value = "  -1.61%  "

# Row 3
$ws.Range("D3").Value = "'2.107.49"
$ws.Range("E3").Value = "  -0.50%  "

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.64%  "

# Row 5
$ws.Range("D5").Value = "'348.80"
$ws.Range("E5").Value = "  +3.56%  "

# Row 6
$ws.Range("E6").Value = "  -0.56%  "

# Row 7
$ws.Range("D7").Value = "'0.5169"
$ws.Range("E7").Value = "  -1.50%  "

# Row 8
$ws.Range("D8").Value = "'0.4456"
$ws.Range("E8").Value = "  -2.33%  "

# Row 9
$ws.Range("D9").Value = "'52.53"
$ws.Range("E9").Value = "  -4.31%  "

# Row 10
$ws.Range("D10").Value = "'0.08990"
$ws.Range("E10").Value = "  -1.28%  "

# Row 11
$ws.Range("E11").Value = "  +0.46%  "

# Row 12
$ws.Range("D12").Value = "'25.74"
$ws.Range("E12").Value = "  +4.50%  "

# Row 13
$ws.Range("D13").Value = "'2.105.27"
$ws.Range("E13").Value = "  -0.74%  "

# Row 14
$ws.Range("D14").Value = "'8.298"
$ws.Range("E14").Value = "  +2.32%  "

# Row 15
$ws.Range("D15").Value = "'6.741"
$ws.Range("E15").Value = "  -1.59%  "

# Row 16
$ws.Range("D16").Value = "'99.45"
$ws.Range("E16").Value = "  +2.43%  "

# Row 17
$ws.Range("D17").Value = "'0.00001152"
$ws.Range("E17").Value = "  -2.09%  "

# Row 18
$ws.Range("D18").Value = "'1.004"
$ws.Range("E18").Value = "  -0.60%  "

# Row 19
$ws.Range("D19").Value = "'20.84"
$ws.Range("E19").Value = "  +7.38%  "

# Row 20
$ws.Range("D20").Value = "'0.06686"
$ws.Range("E20").Value = "  -0.05%  "

# Row 21
$ws.Range("E21").Value = "  -0.51%  "

# Row 22
$ws.Range("D22").Value = "'6.253"
$ws.Range("E22").Value = "  -0.59%  "

# Row 23
$ws.Range("D23").Value = "'30.216.80"
$ws.Range("E23").Value = "  -1.43%  "

# Row 24
$ws.Range("D24").Value = "'12.82"
$ws.Range("E24").Value = "  -0.15%  "

# Row 25
$ws.Range("D25").Value = "'2.347"
$ws.Range("E25").Value = "  -0.60%  "

# Row 26
$ws.Range("D26").Value = "'2.357.42"
$ws.Range("E26").Value = "  -0.29%  "

# Row 27
$ws.Range("D27").Value = "'22.02"
$ws.Range("E27").Value = "  -1.66%  "

# Row 28
$ws.Range("E28").Value = "  +0.11%  "

# Row 29
$ws.Range("D29").Value = "'162.44"
$ws.Range("E29").Value = "  -0.94%  "

# Row 30
$ws.Range("D30").Value = "'133.90"
$ws.Range("E30").Value = "  -0.63%  "

# Row 31
$ws.Range("D31").Value = "'1.183"
$ws.Range("E31").Value = "  -2.37%  "

# Row 32
$ws.Range("D32").Value = "'0.1067"
$ws.Range("E32").Value = "  -0.51%  "

# Row 33
$ws.Range("D33").Value = "'1.641"
$ws.Range("E33").Value = "  -0.10%  "

# Row 34
$ws.Range("D34").Value = "'6.257"
$ws.Range("E34").Value = "  -1.65%  "

# Row 35
$ws.Range("D35").Value = "'3.963"
$ws.Range("E35").Value = "  +0.15%  "

# Row 36
$ws.Range("D36").Value = "'10.31"
$ws.Range("E36").Value = "  -2.28%  "

# Row 37
$ws.Range("D37").Value = "'5.936"
$ws.Range("E37").Value = "  +0.35%  "

# Row 38
$ws.Range("D38").Value = "'0.02579"
$ws.Range("E38").Value = "  -1.76%  "

# Row 39
$ws.Range("E39").Value = "  +0.09%  "

# Row 40
$ws.Range("D40").Value = "'0.2309"
$ws.Range("E40").Value = "  -0.81%  "

# Row 41
$ws.Range("D41").Value = "'12.68"
$ws.Range("E41").Value = "  +0.85%  "

# Row 42
$ws.Range("D42").Value = "'0.6832"
$ws.Range("E42").Value = "  -0.58%  "

# Row 43
$ws.Range("D43").Value = "'1.288"
$ws.Range("E43").Value = "  +2.33%  "

# Row 44
$ws.Range("D44").Value = "'14.28"
$ws.Range("E44").Value = "  -4.12%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.6405"
$ws.Range("E45").Value = "  -0.70%  "

# Row 46
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.306"
$ws.Range("E46").Value = "  -0.37%  "

# Row 47
$ws.Range("D47").Value = "'0.00000000368"
$ws.Range("E47").Value = "  +0.91%  "

# Row 48
$ws.Range("D48").Value = "'3.653"
$ws.Range("E48").Value = "  -1.04%  "

# Row 49
$ws.Range("E49").Value = "  -2.48%  "

# Row 50
$ws.Range("E50").Value = "  -0.58%  "

# Row 51
$ws.Range("D51").Value = "'0.07238"
$ws.Range("E51").Value = "  +0.50%  "
